# Fruta / hortaliza, semanal
# Weekly data refresh: a new week's worth of observations (Primera / Segunda
# quality) is inserted at the top of the "Vega Monumental Concepción - Acelga"
# block, pushing the rest of that block's rows down by two rows (and, in turn,
# appending the two rows that fall off the bottom of the used range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 210; Excel shifts rows 210:253 down to
# 212:255, extending the sheet's dimension from R253 to R255.
$ws.Rows("210:211").Insert()

# New row 210: "Primera" quality observation for the new week.
$ws.Cells.Item(210, 1).Value = 11
$ws.Cells.Item(210, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value = "Bíobío"
$ws.Cells.Item(210, 4).Value = 44694
$ws.Cells.Item(210, 5).Value = 8
$ws.Cells.Item(210, 6).Value = 100112009
$ws.Cells.Item(210, 7).Value = "Acelga"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 200
$ws.Cells.Item(210, 11).Value = 600
$ws.Cells.Item(210, 12).Value = 700
$ws.Cells.Item(210, 13).Value = 650
$ws.Cells.Item(210, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(210, 15).Value = "Región de Ñuble"
$ws.Cells.Item(210, 16).Value = 650
$ws.Cells.Item(210, 17).Value = 1
$ws.Cells.Item(210, 18).Value = "Hortaliza"

# New row 211: "Segunda" quality observation for the same new week.
$ws.Cells.Item(211, 1).Value = 11
$ws.Cells.Item(211, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(211, 3).Value = "Bíobío"
$ws.Cells.Item(211, 4).Value = 44694
$ws.Cells.Item(211, 5).Value = 8
$ws.Cells.Item(211, 6).Value = 100112009
$ws.Cells.Item(211, 7).Value = "Acelga"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Segunda"
$ws.Cells.Item(211, 10).Value = 100
$ws.Cells.Item(211, 11).Value = 500
$ws.Cells.Item(211, 12).Value = 500
$ws.Cells.Item(211, 13).Value = 500
$ws.Cells.Item(211, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(211, 15).Value = "Región de Ñuble"
$ws.Cells.Item(211, 16).Value = 500
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"

Write-Output "Inserted new weekly rows 210-211; sheet now spans to row $($ws.Cells.Item(1,1).Worksheet.UsedRange.Rows.Count)."
